$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue $ws "C2" "21"
Set-TextValue $ws "D2" "42000.00"
Set-TextValue $ws "C3" "86"
Set-TextValue $ws "D3" "305698.00"
Set-TextValue $ws "C5" "141"
Set-TextValue $ws "D5" "382886.40"
Set-TextValue $ws "C6" "412"
Set-TextValue $ws "D6" "1077510.82"
Set-TextValue $ws "C7" "121"
Set-TextValue $ws "D7" "273876.00"
Set-TextValue $ws "C8" "785"
Set-TextValue $ws "D8" "2947058.81"
Set-TextValue $ws "C9" "30"
Set-TextValue $ws "D9" "86600.00"
Set-TextValue $ws "C10" "15"
Set-TextValue $ws "D10" "39500.00"
Set-TextValue $ws "C11" "39"
Set-TextValue $ws "D11" "96177.00"
Set-TextValue $ws "C12" "159"
Set-TextValue $ws "D12" "488316.18"
Set-TextValue $ws "C13" "90"
Set-TextValue $ws "D13" "222800.00"
Set-TextValue $ws "C14" "93"
Set-TextValue $ws "D14" "227788.98"
Set-TextValue $ws "C15" "19"
Set-TextValue $ws "D15" "41593.58"
Set-TextValue $ws "C16" "131"
Set-TextValue $ws "D16" "588717.26"
Set-TextValue $ws "C17" "183"
Set-TextValue $ws "D17" "415089.87"
Set-TextValue $ws "C23" "298"
Set-TextValue $ws "D23" "1161276.10"
Set-TextValue $ws "C35" "167"
Set-TextValue $ws "D35" "473408.00"
Set-TextValue $ws "C37" "376"
Set-TextValue $ws "D37" "1494941.10"
Set-TextValue $ws "C40" "19"
Set-TextValue $ws "D40" "55170.00"
Set-TextValue $ws "C45" "50"
Set-TextValue $ws "D45" "176591.53"
Set-TextValue $ws "C78" "210"
Set-TextValue $ws "D78" "585693.00"
Set-TextValue $ws "C80" "482"
Set-TextValue $ws "D80" "2081939.03"
Set-TextValue $ws "C88" "69"
Set-TextValue $ws "D88" "308536.08"
Set-TextValue $ws "C90" "40"
Set-TextValue $ws "D90" "83555.00"
Set-TextValue $ws "C91" "58"
Set-TextValue $ws "D91" "164457.14"
Set-TextValue $ws "C93" "127"
Set-TextValue $ws "D93" "328107.00"
Set-TextValue $ws "C95" "147"
Set-TextValue $ws "D95" "445497.00"
Set-TextValue $ws "C98" "12"
Set-TextValue $ws "D98" "33000.00"
Set-TextValue $ws "C99" "50"
Set-TextValue $ws "D99" "150000.00"
Set-TextValue $ws "C100" "72"
Set-TextValue $ws "D100" "188591.04"
Set-TextValue $ws "C103" "25"
Set-TextValue $ws "D103" "71330.00"
Set-TextValue $ws "C104" "57"
Set-TextValue $ws "D104" "121492.85"
Set-TextValue $ws "C119" "13"
Set-TextValue $ws "D119" "34000.00"
Set-TextValue $ws "C121" "64"
Set-TextValue $ws "D121" "174877.00"
Set-TextValue $ws "C122" "248"
Set-TextValue $ws "D122" "684008.00"
Set-TextValue $ws "C123" "114"
Set-TextValue $ws "D123" "300781.45"
Set-TextValue $ws "C124" "490"
Set-TextValue $ws "D124" "2190883.06"
Set-TextValue $ws "C127" "31"
Set-TextValue $ws "D127" "72500.00"
Set-TextValue $ws "C128" "88"
Set-TextValue $ws "D128" "270743.68"
Set-TextValue $ws "C129" "43"
Set-TextValue $ws "D129" "158579.76"
Set-TextValue $ws "C130" "56"
Set-TextValue $ws "D130" "161425.82"
Set-TextValue $ws "C132" "87"
Set-TextValue $ws "D132" "385163.75"
Set-TextValue $ws "C133" "121"
Set-TextValue $ws "D133" "304136.44"
Set-TextValue $ws "C140" "2691"
Set-TextValue $ws "D140" "6823017.46"
Set-TextValue $ws "C145" "1055"
Set-TextValue $ws "D145" "2780349.25"
Set-TextValue $ws "C197" "353"
Set-TextValue $ws "D197" "951788.00"
Set-TextValue $ws "C199" "656"
Set-TextValue $ws "D199" "2489424.58"
Set-TextValue $ws "C203" "156"
Set-TextValue $ws "D203" "489133.00"
Set-TextValue $ws "C207" "125"
Set-TextValue $ws "D207" "584388.14"
Set-TextValue $ws "C212" "357"
Set-TextValue $ws "D212" "939423.34"
Set-TextValue $ws "C244" "989"
Set-TextValue $ws "D244" "3594661.86"
